# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Row 21 (CARLOS ENRIQUE VANEGAS CAÑATE / 86677305 / periodo 1609) moves up to
# row 18, and the three "LUIS ENRIQUE JULIO LUQUEZ" / 73100085 rows (previously
# in rows 18-20, sorted 1802/1801/1712) shift down into rows 19-21, now sorted
# ascending by Periodo Mora (1712/1801/1802).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 18 <- Carlos Enrique Vanegas Cañate (was row 21)
$ws.Range("C18").Value = "86677305"
$ws.Range("D18").Value = "CARLOS ENRIQUE VANEGAS CAÑATE"
$ws.Range("E18").Value = "1609"
$ws.Range("F18").Value = 3677
$ws.Range("G18").Value = 689455

# Row 19 <- Luis Enrique Julio Luquez, periodo 1712 (was row 20)
$ws.Range("C19").Value = "73100085"
$ws.Range("D19").Value = "LUIS ENRIQUE JULIO LUQUEZ"
$ws.Range("E19").Value = "1712"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242

# Row 20 <- Luis Enrique Julio Luquez, periodo 1801 (unchanged position)
$ws.Range("C20").Value = "73100085"
$ws.Range("D20").Value = "LUIS ENRIQUE JULIO LUQUEZ"
$ws.Range("E20").Value = "1801"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 781242

# Row 21 <- Luis Enrique Julio Luquez, periodo 1802 (was row 18)
$ws.Range("C21").Value = "73100085"
$ws.Range("D21").Value = "LUIS ENRIQUE JULIO LUQUEZ"
$ws.Range("E21").Value = "1802"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 781242
